$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Financial conditions & employment, US.xlsx",
    "Bitty, SPX & GM2 fitted trends.xlsx",
    "Monetary Expansion and Risk Asset Activity.xlsx",
    "Personal savings U.S (BEA).xlsx",
    "Bitcoin price change is fueled by global monetary growth.xlsx",
    "U.S GDP and GDI.xlsx",
    "Real Incomes U.S (BEA).xlsx",
    "Bitcoin price change is fueled by global monetary growth (with forecast).xlsx",
    "US financial conditions, employment & equities.xlsx",
    "US Equity Indices.xlsx",
    "US Gross Domestic Income and Equity Indices.xlsx",
    "US and Global Monetary Aggregates.xlsx",
    "Bank credit and M2 U.S.xlsx"
)

$row = 2
foreach ($v in $values) {
    $ws.Range("A$row").Value = $v
    $row++
}

$ws.Range("A16").Value = "Monetary Aggregates & Inflation, USA.xlsx"
